$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (first row of the table) gets plain (non-shared) formulas
$ws.Range("B6").Formula = "=`$B`$2"
$ws.Range("C6").Formula = "=B6/((1 + `$B`$3)^A6)"
$ws.Range("D6").Formula = "=SUM(`$C`$6:C6)"

# Rows 7-25 repeat the same relative pattern (mirrors the shared formulas
# Excel itself would have created via fill-down) and each row's "VA Cumulee"
# running total extends the SUM range down to the current row.
for ($r = 7; $r -le 25; $r++) {
    $ws.Range("B$r").Formula = "=`$B`$2"
    $ws.Range("C$r").Formula = "=B$r/((1 + `$B`$3)^A$r)"
    $ws.Range("D$r").Formula = "=SUM(`$C`$6:C$r)"
}

# View changes: zoom to 125% and move the active selection to B11
$excel.ActiveWindow.Zoom = 125
$ws.Range("B11").Select() | Out-Null
